$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 6135
$ws.Range("L3").Value = 6667
$ws.Range("L4").Value = 1641
$ws.Range("L6").Value = 5472
$ws.Range("L7").Value = 20312

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L5").Value = 72
$ws.Range("L8").Value = 1345
$ws.Range("L9").Value = 118
$ws.Range("L16").Value = 46
$ws.Range("L18").Value = 138
$ws.Range("L29").Value = 1138
$ws.Range("L33").Value = 917
$ws.Range("L36").Value = 260
$ws.Range("L37").Value = 779
$ws.Range("L42").Value = 647
$ws.Range("L43").Value = 152
$ws.Range("L44").Value = 138
$ws.Range("L52").Value = 433
$ws.Range("L53").Value = 225
$ws.Range("L63").Value = 61
$ws.Range("L67").Value = 705
$ws.Range("L68").Value = 65
$ws.Range("L76").Value = 308
$ws.Range("L78").Value = 267
$ws.Range("L79").Value = 561
$ws.Range("L83").Value = 445
$ws.Range("L85").Value = 1009
$ws.Range("L86").Value = 131
$ws.Range("L91").Value = 275
$ws.Range("L95").Value = 288
$ws.Range("L96").Value = 226
$ws.Range("L98").Value = 109
$ws.Range("L99").Value = 354
$ws.Range("L100").Value = 39
$ws.Range("L101").Value = 20312

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L2").Value = 72
$ws.Range("L7").Value = 226

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 306
$ws.Range("L3").Value = 419
$ws.Range("L7").Value = 1009

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L6").Value = 123
$ws.Range("L7").Value = 433

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L6").Value = 76
$ws.Range("L7").Value = 225

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 404
$ws.Range("L3").Value = 477
$ws.Range("L7").Value = 1345

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L6").Value = 99
$ws.Range("L7").Value = 445

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 248
$ws.Range("L7").Value = 917

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L3").Value = 93
$ws.Range("L7").Value = 288

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 235
$ws.Range("L7").Value = 779

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L3").Value = 144
$ws.Range("L7").Value = 354

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 200
$ws.Range("L6").Value = 163
$ws.Range("L7").Value = 705

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L3").Value = 440
$ws.Range("L6").Value = 279
$ws.Range("L7").Value = 1138

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L6").Value = 34
$ws.Range("L7").Value = 138

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L6").Value = 137
$ws.Range("L7").Value = 308

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 177
$ws.Range("L6").Value = 180
$ws.Range("L7").Value = 647

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L6").Value = 76
$ws.Range("L7").Value = 267

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L2").Value = 94
$ws.Range("L3").Value = 125
$ws.Range("L4").Value = 15
$ws.Range("L7").Value = 275

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L3").Value = 180
$ws.Range("L7").Value = 561

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("L2").Value = 49
$ws.Range("L7").Value = 138

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L3").Value = 84
$ws.Range("L7").Value = 260

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("L6").Value = 24
$ws.Range("L7").Value = 39

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("L6").Value = 49
$ws.Range("L7").Value = 109

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("L6").Value = 28
$ws.Range("L7").Value = 118

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("L6").Value = 32
$ws.Range("L7").Value = 72

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("L4").Value = 70
$ws.Range("L7").Value = 131

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L2").Value = 70
$ws.Range("L6").Value = 61

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("L2").Value = 22
$ws.Range("L7").Value = 65

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L6").Value = 47
$ws.Range("L7").Value = 152

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("L6").Value = 30
$ws.Range("L7").Value = 46
